$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$rows = 8,9,10,11,13,14

foreach ($r in $rows) {
    # Priority column (E) changes from blank to "ht" on both language sheets
    $zh.Cells.Item($r, 5).Value = "ht"
    $de.Cells.Item($r, 5).Value = "ht"

    # Latest Handoff Datetime (column H) on zh-cn sheet
    $zh.Cells.Item($r, 8).Value = "2016-08-25 02:21:25"

    # Latest Handoff Datetime (column H) on de-de sheet, shared text with Overview column G
    $de.Cells.Item($r, 8).Value = "2016-08-25 02:21:30"

    # Latest HO Xliff Generate Date (column G) on Overview sheet
    $overview.Cells.Item($r, 7).Value = "2016-08-25 02:21:30"
}
